$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-31 Wednesday", "2024-02-01 Thursday"),
    @("13×81=1053", "78×52=4056"),
    @("92×38=3496", "32×73=2336"),
    @("62×81=5022", "54×77=4158"),
    @("92×83=7636", "50×21=1050"),
    @("42×56=2352", "54×13=702"),
    @("85×56=4760", "96×32=3072"),
    @("44×18=792", "18×31=558"),
    @("43×71=3053", "21×35=735"),
    @("32×94=3008", "98×38=3724"),
    @("41×63=2583", "40×63=2520"),
    @("17×85=1445", "11×60=660"),
    @("51×87=4437", "27×42=1134"),
    @("86×14=1204", "48×67=3216"),
    @("99×78=7722", "64×85=5440"),
    @("29×56=1624", "72×67=4824"),
    @("50×86=4300", "31×60=1860"),
    @("44×61=2684", "81×45=3645"),
    @("35×44=1540", "96×69=6624"),
    @("65×83=5395", "66×75=4950"),
    @("19×86=1634", "33×78=2574"),
    @("19×84=1596", "33×28=924"),
    @("56×63=3528", "64×71=4544"),
    @("26×32=832", "60×30=1800"),
    @("43×85=3655", "39×88=3432"),
    @("17×48=816", "60×73=4380")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
